$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# "Marking" row: Right marks per correct answer (B11): 3 -> 5
$ws.Range("B11").Value = 5

# "Total" row: total marks obtained (B12): 15 -> 25
$ws.Range("B12").Value = 25

# "Total" row: correct/total max text (E12): "14/84" -> "25/140"
$ws.Range("E12").Value = "25/140"
